# Apply text replacements per the diff: update the date line and the
# two-digit multiplication problems scattered through the table cells.

$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-12 Saturday", "2023-08-13 Sunday"),
    @("67×29=", "67×52="),
    @("79×41=", "25×97="),
    @("64×73=", "84×93="),
    @("14×43=", "11×49="),
    @("15×95=", "48×79="),
    @("28×34=", "81×47="),
    @("90×17=", "46×52="),
    @("93×23=", "37×18="),
    @("48×66=", "43×11="),
    @("42×40=", "13×19="),
    @("78×23=", "91×40="),
    @("51×50=", "80×74="),
    @("46×13=", "36×34="),
    @("85×35=", "23×64="),
    @("96×87=", "96×63="),
    @("26×15=", "97×60="),
    @("95×17=", "33×84="),
    @("84×64=", "71×81="),
    @("31×19=", "29×94="),
    @("94×18=", "17×82="),
    @("18×62=", "64×20="),
    @("12×34=", "37×45="),
    @("12×30=", "85×69="),
    @("86×14=", "99×85="),
    @("50×26=", "46×11=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
